$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: new literature entry (Afrouzi et al 2020, working memory)
# Values are entered in the order A18, C18, B18, D18 to match the
# shared-string insertion order of the original edit.
$ws.Cells.Item(18, 1).Value = 'Afrouzi et al WP 2020 working memory'
$ws.Cells.Item(18, 3).Value = 'cgain learning and diagnostic overreact too little. A main feature in the data is the variation of overreaction across different settings, and that overreaction appears to be stronger when the forecast horizon is longer (see Bouchaud et al. (2019) and Bordalo et al. (2019) for evidence from analyst earnings forecasts, as well as Brooks, Katz and Lustig (2018), Wang (2019), and d’Arienzo (2020) for evidence from interest rate forecasts). '
$ws.Cells.Item(18, 2).Value = 'Experiment of expectations, documenting new evidence. They develop a working memory model, where ppl estimate long-run means of the process subject to a cost of utilizing past information. '
$ws.Cells.Item(18, 4).Value = 'Should be my standard response to the diagnostic E question.'

# Match the yellow highlight + wrap-text formatting used by the other
# "important" rows (e.g. row 4, row 10).
$newRow = $ws.Range("A18:D18")
$newRow.Interior.Color = 65535
$newRow.WrapText = $true
$ws.Rows.Item(18).RowHeight = 75

# Update the view: scroll/selection moved to column B/ cell C18.
$ws.Cells.Item(18, 3).Select()

